$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value   = -20.82769999999998
$ws.Range("A7").Value   = -19.33949999999998
$ws.Range("A16").Value  = -21.67169999999999
$ws.Range("A28").Value  = -22.0951
$ws.Range("A29").Value  = -21.44339999999998
$ws.Range("A32").Value  = -21.15099999999999
$ws.Range("A40").Value  = -20.1697
$ws.Range("A52").Value  = -22.1111
$ws.Range("A57").Value  = -22.3426
$ws.Range("A66").Value  = -21.48639999999999
$ws.Range("A100").Value = -21.8332
